# Rename the existing sheet to ID_0001 and populate it with the
# username/password table.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ID_0001"

$ws1.Range("A1").Value = "USERNAME"
$ws1.Range("B1").Value = "PASSWORD"
$ws1.Range("A2").Value = "automation.devmrkolv@gmail.com"
$ws1.Range("B2").Value = "`$chlUe13elKiNd"
$ws1.Range("A3").Value = "mmm"
$ws1.Range("B3").Value = "m"
$ws1.Range("A4").Value = "z."
$ws1.Range("B4").Value = "z"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "mailto:automation.devmrkolv@gmail.com")
$ws1.Range("A2").Style = "Normal"

$ws1.Columns.Item(1).AutoFit() | Out-Null
$ws1.Columns.Item(2).AutoFit() | Out-Null

$ws1.Range("A5").Select() | Out-Null

# Add a second sheet that mirrors ID_0001 via formulas.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "ID_0002"

$ws2.Range("A1").Formula = "=ID_0001!A1"
$ws2.Range("B1").Formula = "=ID_0001!B1"
$ws2.Range("A2").Formula = "=ID_0001!A2"
$ws2.Range("B2").Formula = "=ID_0001!B2"
$ws2.Range("A3").Formula = "=ID_0001!A3"
$ws2.Range("B3").Formula = "=ID_0001!B3"
$ws2.Range("A4").Formula = "=ID_0001!A4"
$ws2.Range("B4").Formula = "=ID_0001!B4"

$ws2.PageSetup.LeftMargin = 36.850393728
$ws2.PageSetup.RightMargin = 36.850393728
$ws2.PageSetup.TopMargin = 56.692913399999995
$ws2.PageSetup.BottomMargin = 56.692913399999995
$ws2.PageSetup.HeaderMargin = 22.67716464
$ws2.PageSetup.FooterMargin = 22.67716464

$ws2.Range("G9").Select() | Out-Null
$ws2.Activate()
